# Update test data d1.xlsx: fix expected output for invalid input
# (home price <= 0 should show a validation message instead of a $0.00 payment)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# E5 held the stale "Monthly Pay:   $0.00" text for the negative-home-price
# test row; it should now hold the validation message that the app actually
# returns for a non-positive home price.
$ws.Range("E5").Value = "Please provide a positive home price value."

# Rename the built-in "Hyperlink" cell style to "Link" (cosmetic rename).
$wb.Styles.Item("Hyperlink").Name = "Link"

# Move/save the sheet's active selection from E13 to E10.
$ws.Range("E10").Select()

# Persist the workbook window's last-used size.
$excel.ActiveWindow.Width = 13800
$excel.ActiveWindow.Height = 12420
